# Update attendance/view numbers (column F) across sheets as per the
# upstream data refresh (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F20").Value = 1856
$ws.Range("F21").Value = 1387
$ws.Range("F30").Value = 4723
$ws.Range("F31").Value = 2280
$ws.Range("F32").Value = 3924
$ws.Range("F34").Value = 168
$ws.Range("F48").Value = 14
$ws.Range("F49").Value = 165

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 97

# Sheet: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 2914
$ws.Range("F15").Value = 424

# Sheet: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 2914
$ws.Range("F19").Value = 1856
$ws.Range("F20").Value = 1387
$ws.Range("F26").Value = 97
$ws.Range("F34").Value = 4723
$ws.Range("F35").Value = 2280
$ws.Range("F36").Value = 3924
$ws.Range("F38").Value = 168
